$d = $word.ActiveDocument

# --- Change 1: merge "Professor " + "Shoarinejad" runs (and drop the
#     spell-check proofErr markers around "Shoarinejad") into a single run.
#     We locate the paragraph by index, insert a brand new paragraph right
#     before it (new paragraphs never carry stray <w:proofErr/> markers),
#     give the new paragraph the merged text, then delete the old one.
$profIdx = 0
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Professor Shoarinejad`r") {
        $profIdx = $i
        break
    }
    $i = $i + 1
}

if ($profIdx -gt 0) {
    $prevPara = $d.Paragraphs($profIdx - 1)
    $prevPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs($profIdx)
    $newPara.Range.Text = "Professor Shoarinejad"
    # the original paragraph has now shifted one slot further down
    $oldPara = $d.Paragraphs($profIdx + 1)
    $oldPara.Range.Delete()
}

# --- Change 2: add a new bullet after the "static IP address" paragraph
#     describing the Tutorial 2 / MQTT work.
$ipPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*static IP address*") {
        $ipPara = $p
        break
    }
}
if ($ipPara -ne $null) {
    $ipPara.Range.InsertParagraphAfter()
    $newPara2 = $ipPara.Next()
    $newPara2.Range.Text = "For Tutorial 2, I went through and did the TCP/IP connection tutorial (Part 2) and ran into some issues. I talked to a TA and was able to get past that part and pushed the code onto Github."
}
